$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Discretionary")
$ws.Range("B2").Value = -38915.87
$ws.Range("C2").Value = -42885.77
$ws.Range("D2").Value = -38293.06
$ws.Range("E2").Value = -31954.07
$ws.Range("F2").Value = -28973.95
$ws.Range("G2").Value = -181022.72
$ws.Range("B3").Value = -11866.61
$ws.Range("C3").Value = -20344.57
$ws.Range("D3").Value = -15438.74
$ws.Range("E3").Value = -15008.15
$ws.Range("F3").Value = -18061.25
$ws.Range("G3").Value = -80719.32
$ws.Range("B4").Value = -42160.42
$ws.Range("C4").Value = -40308.66
$ws.Range("D4").Value = -31167.05
$ws.Range("E4").Value = -37792.51
$ws.Range("F4").Value = -39080.4
$ws.Range("G4").Value = -190509.04
$ws.Range("B5").Value = -11891.18
$ws.Range("C5").Value = -20277.56
$ws.Range("D5").Value = -16132.08
$ws.Range("E5").Value = -18620.19
$ws.Range("F5").Value = -13311.51
$ws.Range("G5").Value = -80232.52
$ws.Range("B6").Value = -57936.9
$ws.Range("C6").Value = -44474.84
$ws.Range("D6").Value = -54028.1
$ws.Range("E6").Value = -73444.08
$ws.Range("F6").Value = -49862.93
$ws.Range("G6").Value = -279746.85
$ws.Range("B7").Value = -28564.63
$ws.Range("C7").Value = -29994.23
$ws.Range("D7").Value = -27644.43
$ws.Range("E7").Value = -29930.08
$ws.Range("F7").Value = -24965.85
$ws.Range("G7").Value = -141099.22
$ws.Range("B8").Value = -191335.61
$ws.Range("C8").Value = -198285.63
$ws.Range("D8").Value = -182703.46
$ws.Range("E8").Value = -206749.08
$ws.Range("F8").Value = -174255.89
$ws.Range("G8").Value = -953329.67

$ws = $wb.Worksheets.Item("Essential")
$ws.Range("B2").Value = -36353.05
$ws.Range("C2").Value = -30457.38
$ws.Range("D2").Value = -38819.41
$ws.Range("E2").Value = -38946.09
$ws.Range("F2").Value = -34430.07
$ws.Range("G2").Value = -179006
$ws.Range("B3").Value = -18135.67
$ws.Range("C3").Value = -16482.58
$ws.Range("D3").Value = -15629.56
$ws.Range("E3").Value = -19907.67
$ws.Range("F3").Value = -12846.96
$ws.Range("G3").Value = -83002.44
$ws.Range("B4").Value = -54488.72
$ws.Range("C4").Value = -46939.96
$ws.Range("D4").Value = -54448.97
$ws.Range("E4").Value = -58853.76
$ws.Range("F4").Value = -47277.03
$ws.Range("G4").Value = -262008.44

$ws = $wb.Worksheets.Item("Income")
$ws.Range("B2").Value = 57327.05
$ws.Range("C2").Value = 69089.91
$ws.Range("D2").Value = 46505.97
$ws.Range("E2").Value = 68042.99
$ws.Range("F2").Value = 97647.47
$ws.Range("G2").Value = 338613.39
$ws.Range("B3").Value = 452802.19
$ws.Range("C3").Value = 439948.04
$ws.Range("D3").Value = 521104.54
$ws.Range("E3").Value = 400785.61
$ws.Range("F3").Value = 439761.29
$ws.Range("G3").Value = 2254401.67
$ws.Range("B4").Value = 81716.1
$ws.Range("C4").Value = 76255.21
$ws.Range("D4").Value = 76187.9
$ws.Range("E4").Value = 68580.99
$ws.Range("F4").Value = 110584.34
$ws.Range("G4").Value = 413324.54
$ws.Range("B5").Value = 591845.34
$ws.Range("C5").Value = 585293.16
$ws.Range("D5").Value = 643798.41
$ws.Range("E5").Value = 537409.59
$ws.Range("F5").Value = 647993.1
$ws.Range("G5").Value = 3006339.6

$ws = $wb.Worksheets.Item("Transfer")
$ws.Range("B2").Value = 196582.08
$ws.Range("C2").Value = 188143.44
$ws.Range("D2").Value = 185434.74
$ws.Range("E2").Value = 205682.74
$ws.Range("F2").Value = 189920.85
$ws.Range("G2").Value = 965763.85
$ws.Range("B3").Value = -196582.08
$ws.Range("C3").Value = -188143.44
$ws.Range("D3").Value = -185434.74
$ws.Range("E3").Value = -205682.74
$ws.Range("F3").Value = -189920.85
$ws.Range("G3").Value = -965763.85
